$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5
$ws.Range("E5").Value = 21
$ws.Range("F5").Value = 11
$ws.Range("H5").Value = 11

# Row 6
$ws.Range("E6").Value = 42
$ws.Range("F6").Value = 22
$ws.Range("H6").Value = 22

# Row 8
$ws.Range("E8").Value = 30

# Row 10
$ws.Range("E10").Value = 18

# Row 11
$ws.Range("E11").Value = 10

# Row 12
$ws.Range("E12").Value = 20

# Row 16
$ws.Range("E16").Value = 254
$ws.Range("F16").Value = 69
$ws.Range("H16").Value = 69

# Row 18
$ws.Range("E18").Value = 71
